$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2694.1667
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 224
$ws.Range("I5").Value = 224
$ws.Range("K5").Value = 224
$ws.Range("M5").Value = -112
$ws.Range("H14").Value = 2383.8333
$ws.Range("I14").Value = 2563.3635
$ws.Range("J14").Value = 409
$ws.Range("K14").Value = 2563.3635
$ws.Range("L14").Value = 409
$ws.Range("M14").Value = -2388.3635
$ws.Range("N14").Value = -759
$ws.Range("H15").Value = 375
$ws.Range("I15").Value = 250
$ws.Range("J15").Value = 500
$ws.Range("K15").Value = 250
$ws.Range("L15").Value = 500
$ws.Range("M15").Value = 100
$ws.Range("N15").Value = -1200
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H22").Value = 8888.111000000001
$ws.Range("I22").Value = 16004
$ws.Range("J22").Value = 3195.4
$ws.Range("K22").Value = 16004
$ws.Range("L22").Value = 3195.4
$ws.Range("M22").Value = -15705
$ws.Range("N22").Value = -3793.4
$ws.Range("H36").Value = 1649.875
$ws.Range("I36").Value = 1649.875
$ws.Range("K36").Value = 1649.875
$ws.Range("M36").Value = -1303.875
$ws.Range("H50").Value = 1001
$ws.Range("I50").Value = 725
$ws.Range("J50").Value = 1553
$ws.Range("K50").Value = 725
$ws.Range("L50").Value = 1553
$ws.Range("M50").Value = -11
$ws.Range("N50").Value = -2981
$ws.Range("H58").Value = 99998
$ws.Range("J58").Value = 99998
$ws.Range("L58").Value = 99998
$ws.Range("N58").Value = -100858
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 224
$ws.Range("I4").Value = 224
$ws.Range("K4").Value = 224
$ws.Range("M4").Value = -109
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 5001750
$ws.Range("I3").Value = 5001750
$ws.Range("K3").Value = 5001750
$ws.Range("M3").Value = -5001637
$ws.Range("H16").Value = 1209.2
$ws.Range("I16").Value = 1170.4445
$ws.Range("J16").Value = 1267.3334
$ws.Range("K16").Value = 1170.4445
$ws.Range("L16").Value = 1267.3334
$ws.Range("M16").Value = -883.4445000000001
$ws.Range("N16").Value = -1841.3334
$ws.Range("H105").Value = 816.7857
$ws.Range("I105").Value = 794.5833
$ws.Range("K105").Value = 794.5833
$ws.Range("M105").Value = 952.4167
$ws.Range("H113").Value = 1209.2
$ws.Range("I113").Value = 1170.4445
$ws.Range("J113").Value = 1267.3334
$ws.Range("K113").Value = 1170.4445
$ws.Range("L113").Value = 1267.3334
$ws.Range("M113").Value = 999.5554999999999
$ws.Range("N113").Value = -5607.3334
$ws.Range("H141").Value = 85843
$ws.Range("J141").Value = 85843
$ws.Range("L141").Value = 85843
$ws.Range("N141").Value = -96203
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 166821.67
$ws.Range("I4").Value = 250125.25
$ws.Range("J4").Value = 214.5
$ws.Range("K4").Value = 750375.75
$ws.Range("L4").Value = 643.5
$ws.Range("M4").Value = -750263.75
$ws.Range("N4").Value = -867.5
$ws.Range("H6").Value = 35
$ws.Range("I6").Value = 35
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 105
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 8
$ws.Range("N6").ClearContents()
$ws.Range("H11").Value = 140.36363
$ws.Range("I11").Value = 48.5
$ws.Range("K11").Value = 145.5
$ws.Range("M11").Value = -5.5
$ws.Range("H17").Value = 3400
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 3400
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 10200
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -10538
$ws.Range("H69").Value = 1225
$ws.Range("J69").Value = 1225
$ws.Range("L69").Value = 3675
$ws.Range("N69").Value = -5297
$ws.Range("H72").Value = 1225
$ws.Range("J72").Value = 1225
$ws.Range("L72").Value = 11025
$ws.Range("N72").Value = -19137
$ws.Range("H74").Value = 17500
$ws.Range("I74").Value = 10000
$ws.Range("J74").Value = 25000
$ws.Range("K74").Value = 30000
$ws.Range("L74").Value = 75000
$ws.Range("M74").Value = -28939
$ws.Range("N74").Value = -77122
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("H77").Value = 17500
$ws.Range("I77").Value = 10000
$ws.Range("J77").Value = 25000
$ws.Range("K77").Value = 90000
$ws.Range("L77").Value = 225000
$ws.Range("M77").Value = -84696
$ws.Range("N77").Value = -235608
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 208.1
$ws.Range("I107").Value = 91.833336
$ws.Range("J107").Value = 382.5
$ws.Range("K107").Value = 91.833336
$ws.Range("L107").Value = 382.5
$ws.Range("M107").Value = 1828.166664
$ws.Range("N107").Value = -4222.5
$ws.Range("H113").Value = 8938.25
$ws.Range("I113").Value = 1499
$ws.Range("K113").Value = 1499
$ws.Range("M113").Value = 671
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3509
$ws.Range("I82").Value = 1359
$ws.Range("J82").Value = 4584
$ws.Range("K82").Value = 1359
$ws.Range("L82").Value = 4584
$ws.Range("M82").Value = -998
$ws.Range("N82").Value = -5306
$ws.Range("H85").Value = 3509
$ws.Range("I85").Value = 1359
$ws.Range("J85").Value = 4584
$ws.Range("K85").Value = 1359
$ws.Range("L85").Value = 4584
$ws.Range("M85").Value = -111
$ws.Range("N85").Value = -7080
$ws.Range("H93").Value = 1216.8572
$ws.Range("I93").Value = 1216.8572
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 1216.8572
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 31.14280000000008
$ws.Range("N93").ClearContents()
$ws.Range("H127").Value = 69500.5
$ws.Range("J127").Value = 69500.5
$ws.Range("L127").Value = 69500.5
$ws.Range("N127").Value = -79420.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 5000
$ws.Range("I3").Value = 5000
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 5000
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -4886
$ws.Range("N3").ClearContents()
